$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 19 ("Typical Lease Term" / 24 / "months") entirely - everything
# below shifts up by one row.
$ws.Rows.Item(19).Delete()

# The cell that used to read "Monthly Lease" (now at A22 after the shift)
# should be renamed to "Yearly Lease".
$ws.Range("A22").Value = "Yearly Lease"

# Update the active selection to match the saved workbook view.
$ws.Range("A23").Select()
